$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9937165061569999
$ws.Range("D2").Value = 0.02883949465213931
$ws.Range("E2").Value = 0.3059238889258449
$ws.Range("F2").Value = 0.7565922461509373
$ws.Range("G2").Value = 0.002430152510010241
$ws.Range("K2").Value = 0.3911746362539077
$ws.Range("L2").Value = 0.1396917977880463
$ws.Range("N2").Value = 1.837630429286745
$ws.Range("O2").Value = 2.607734100147212
$ws.Range("B3").Value = 0.9679792446818567
$ws.Range("D3").Value = 0.02712728148494392
$ws.Range("E3").Value = 0.3083394191694939
$ws.Range("F3").Value = 0.7524161873701729
$ws.Range("G3").Value = 0.002432617965832573
$ws.Range("K3").Value = 0.3415326572455228
$ws.Range("L3").Value = 0.1291995543510325
$ws.Range("N3").Value = 1.856865533837679
$ws.Range("O3").Value = 2.607272311714127
$ws.Range("B4").Value = 0.9526466160103269
$ws.Range("D4").Value = 0.02606572378602579
$ws.Range("E4").Value = 0.3099375296855698
$ws.Range("F4").Value = 0.7503063015120475
$ws.Range("G4").Value = 0.002434213931948306
$ws.Range("K4").Value = 0.3109323911908746
$ws.Range("L4").Value = 0.1228123562085131
$ws.Range("N4").Value = 1.869284042096816
$ws.Range("O4").Value = 2.608475284484712
$ws.Range("B5").Value = 0.9465171930580993
$ws.Range("D5").Value = 0.02563057261146184
$ws.Range("E5").Value = 0.3106177149966793
$ws.Range("F5").Value = 0.7495607584654564
$ws.Range("G5").Value = 0.002434885023632404
$ws.Range("K5").Value = 0.2984329840475368
$ws.Range("L5").Value = 0.1202234666762791
$ws.Range("N5").Value = 1.874497603653925
$ws.Range("O5").Value = 2.609339317128416
$ws.Range("B6").Value = 0.9455065945334127
$ws.Range("D6").Value = 0.02555816219761908
$ws.Range("E6").Value = 0.3107324082750331
$ws.Range("F6").Value = 0.749443863148926
$ws.Range("G6").Value = 0.002434997711255324
$ws.Range("K6").Value = 0.2963556998396086
$ws.Range("L6").Value = 0.1197944286137869
$ws.Range("N6").Value = 1.875372547908822
$ws.Range("O6").Value = 2.609505366891995
$ws.Range("B7").Value = 0.952563470969352
$ws.Range("D7").Value = 0.02605986551576933
$ws.Range("E7").Value = 0.3099465856714012
$ws.Range("F7").Value = 0.7502957842113176
$ws.Range("G7").Value = 0.002434222898545008
$ws.Range("K7").Value = 0.3107639386710446
$ws.Range("L7").Value = 0.1227773849476961
$ws.Range("N7").Value = 1.869353734832313
$ws.Range("O7").Value = 2.608485423564844
$ws.Range("B8").Value = 0.9847450218603626
$ws.Range("D8").Value = 0.02825126294955993
$ws.Range("E8").Value = 0.3067329228366305
$ws.Range("F8").Value = 0.7550580814646466
$ws.Range("G8").Value = 0.002430985582394892
$ws.Range("K8").Value = 0.374083308833832
$ws.Range("L8").Value = 0.1360626935527733
$ws.Range("N8").Value = 1.844136496585486
$ws.Range("O8").Value = 2.607266351533013
$ws.Range("B9").Value = 1.051564598015176
$ws.Range("D9").Value = 0.03246653854565551
$ws.Range("E9").Value = 0.3013416486353524
$ws.Range("F9").Value = 0.7680012245628376
$ws.Range("G9").Value = 0.002425286383151063
$ws.Range("K9").Value = 0.4972806117490904
$ws.Range("L9").Value = 0.1625494250464499
$ws.Range("N9").Value = 1.799509062023088
$ws.Range("O9").Value = 2.616673570146446
$ws.Range("B10").Value = 1.102900580002682
$ws.Range("D10").Value = 0.03551280235468113
$ws.Range("E10").Value = 0.297933838549957
$ws.Range("F10").Value = 0.7797104259592231
$ws.Range("G10").Value = 0.002421491030728275
$ws.Range("K10").Value = 0.5871821032734488
$ws.Range("L10").Value = 0.1822723227745513
$ws.Range("N10").Value = 1.769658930352516
$ws.Range("O10").Value = 2.63078683809573
$ws.Range("B11").Value = 1.126738159753302
$ws.Range("D11").Value = 0.03688748805136299
$ws.Range("E11").Value = 0.2965032225593216
$ws.Range("F11").Value = 0.7855156264986078
$ws.Range("G11").Value = 0.002419848683946846
$ws.Range("K11").Value = 0.6279442576193048
$ws.Range("L11").Value = 0.1913016768281608
$ws.Range("N11").Value = 1.75671660395042
$ws.Range("O11").Value = 2.638773524871283
$ws.Range("B12").Value = 1.135834068209675
$ws.Range("D12").Value = 0.03740643577150848
$ws.Range("E12").Value = 0.2959786533819226
$ws.Range("F12").Value = 0.7877827386694776
$ws.Range("G12").Value = 0.002419238811902006
$ws.Range("K12").Value = 0.6433599869158968
$ws.Range("L12").Value = 0.1947290307488601
$ws.Range("N12").Value = 1.75190722633387
$ws.Range("O12").Value = 2.642023228633121
$ws.Range("B13").Value = 1.133872035545266
$ws.Range("D13").Value = 0.03729474333592719
$ws.Range("E13").Value = 0.2960908653369465
$ws.Range("F13").Value = 0.7872914154981743
$ws.Range("G13").Value = 0.002419369623685671
$ws.Range("K13").Value = 0.6400408346658253
$ws.Range("L13").Value = 0.1939905287360233
$ws.Range("N13").Value = 1.75293893705073
$ws.Range("O13").Value = 2.641313324740821
$ws.Range("B14").Value = 1.127485103203412
$ws.Range("D14").Value = 0.03693021471261204
$ws.Range("E14").Value = 0.2964597218939655
$ws.Range("F14").Value = 0.7857007642543294
$ws.Range("G14").Value = 0.002419798268158321
$ws.Range("K14").Value = 0.629212923922978
$ws.Range("L14").Value = 0.1915834849078948
$ws.Range("N14").Value = 1.756319097445187
$ws.Range("O14").Value = 2.639036363897191
$ws.Range("B15").Value = 1.123581909775709
$ws.Range("D15").Value = 0.03670671927152824
$ws.Range("E15").Value = 0.2966878929760011
$ws.Range("F15").Value = 0.7847354048366242
$ws.Range("G15").Value = 0.002420062392928246
$ws.Range("K15").Value = 0.6225778894997518
$ws.Range("L15").Value = 0.1901101570849022
$ws.Range("N15").Value = 1.758401473946798
$ws.Range("O15").Value = 2.637671004527476
$ws.Range("B16").Value = 1.101352442324156
$ws.Range("D16").Value = 0.03542273855929068
$ws.Range("E16").Value = 0.298029737357334
$ws.Range("F16").Value = 0.7793406721929159
$ws.Range("G16").Value = 0.002421600051256037
$ws.Range("K16").Value = 0.584515436737604
$ws.Range("L16").Value = 0.1816833772817148
$ws.Range("N16").Value = 1.770517557127356
$ws.Range("O16").Value = 2.630296416241407
$ws.Range("B17").Value = 1.087839107771657
$ws.Range("D17").Value = 0.03463220481263818
$ws.Range("E17").Value = 0.2988835326704109
$ws.Range("F17").Value = 0.776153754240994
$ws.Range("G17").Value = 0.00242256487507068
$ws.Range("K17").Value = 0.5611304297209472
$ws.Range("L17").Value = 0.1765284227300157
$ws.Range("N17").Value = 1.778113478175429
$ws.Range("O17").Value = 2.626173603787862
$ws.Range("B18").Value = 1.080112239310893
$ws.Range("D18").Value = 0.03417646967964316
$ws.Range("E18").Value = 0.2993858733442973
$ws.Range("F18").Value = 0.7743657749709314
$ws.Range("G18").Value = 0.002423127742566948
$ws.Range("K18").Value = 0.5476673742026321
$ws.Range("L18").Value = 0.1735688298654452
$ws.Range("N18").Value = 1.782542372146404
$ws.Range("O18").Value = 2.623949713144441
$ws.Range("B19").Value = 1.077503910538951
$ws.Range("D19").Value = 0.03402198750725915
$ws.Range("E19").Value = 0.299557892016697
$ws.Range("F19").Value = 0.7737681340804983
$ws.Range("G19").Value = 0.002423319683122333
$ws.Range("K19").Value = 0.5431068683436706
$ws.Range("L19").Value = 0.1725676934056111
$ws.Range("N19").Value = 1.784052209538529
$ws.Range("O19").Value = 2.623222063463459
$ws.Range("B20").Value = 1.089272905274157
$ws.Range("D20").Value = 0.03471646639195569
$ws.Range("E20").Value = 0.2987914795053754
$ws.Range("F20").Value = 0.7764883445411925
$ws.Range("G20").Value = 0.002422461347851623
$ws.Range("K20").Value = 0.5636211144435777
$ws.Range("L20").Value = 0.1770766182729631
$ws.Range("N20").Value = 1.777298677439787
$ws.Range("O20").Value = 2.626597224075113
$ws.Range("B21").Value = 1.129359228140373
$ws.Range("D21").Value = 0.03703732964712003
$ws.Range("E21").Value = 0.2963509138919616
$ws.Range("F21").Value = 0.7861661097759765
$ws.Range("G21").Value = 0.002419672038439088
$ws.Range("K21").Value = 0.6323938923344485
$ws.Range("L21").Value = 0.1922902721463515
$ws.Range("N21").Value = 1.75532377596684
$ws.Range("O21").Value = 2.639699046938887
$ws.Range("B22").Value = 1.155960670433785
$ws.Range("D22").Value = 0.0385447182608516
$ws.Range("E22").Value = 0.2948559545794396
$ws.Range("F22").Value = 0.7928921499170656
$ws.Range("G22").Value = 0.002417919269228745
$ws.Range("K22").Value = 0.6772239084157832
$ws.Range("L22").Value = 0.2022806324489039
$ws.Range("N22").Value = 1.741495851183423
$ws.Range("O22").Value = 2.649575202362314
$ws.Range("B23").Value = 1.141726313095944
$ws.Range("D23").Value = 0.03774106713490966
$ws.Range("E23").Value = 0.2956446927801348
$ws.Range("F23").Value = 0.7892656445446278
$ws.Range("G23").Value = 0.002418848349540759
$ws.Range("K23").Value = 0.6533082110464363
$ws.Range("L23").Value = 0.1969442915197419
$ws.Range("N23").Value = 1.748827201329661
$ws.Range("O23").Value = 2.644183920945864
$ws.Range("B24").Value = 1.088624554159253
$ws.Range("D24").Value = 0.0346783756136233
$ws.Range("E24").Value = 0.2988330609734469
$ws.Range("F24").Value = 0.7763369385153993
$ws.Range("G24").Value = 0.002422508127004007
$ws.Range("K24").Value = 0.5624951340902555
$ws.Range("L24").Value = 0.1768287664172448
$ws.Range("N24").Value = 1.777666856045674
$ws.Range("O24").Value = 2.626405249413978
$ws.Range("B25").Value = 1.033092745655381
$ws.Range("D25").Value = 0.03133504879142635
$ws.Range("E25").Value = 0.3027028442699411
$ws.Range("F25").Value = 0.7641136249657663
$ws.Range("G25").Value = 0.002426759074490533
$ws.Range("K25").Value = 0.4640584914843373
$ws.Range("L25").Value = 0.1553377387450467
$ws.Range("N25").Value = 1.811065887261281
$ws.Range("O25").Value = 2.609339317128416
